# NatmiData LR-pairs update: Efna1-Epha2, OldD0 TPM refresh.
# The underlying average-expression inputs for the "ECs" cluster were
# recomputed with the new TPM values (Ligand avg expr for ECs as sender,
# Receptor avg expr for ECs as target). Every other column in this sheet
# (total expression, specificities, edge weights) is purely derived from
# those two per-cluster averages, so this script writes the fully
# recomputed values for every affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> ECs -------------------------------------------------
$ws.Range("G2").Value = 23.630375
$ws.Range("H2").Value = 70.89112499999999
$ws.Range("I2").Value = 0.9002398112414131
$ws.Range("J2").Value = 0.9002398112414129
$ws.Range("M2").Value = 10.858287
$ws.Range("N2").Value = 32.574861
$ws.Range("O2").Value = 0.5084025289165609
$ws.Range("P2").Value = 0.508402528916561
$ws.Range("Q2").Value = 256.5853936676249
$ws.Range("R2").Value = 2309.268543008624
$ws.Range("S2").Value = 0.4576841966665018
$ws.Range("T2").Value = 0.4576841966665019

# --- Row 3: ECs -> FAPs -------------------------------------------------
$ws.Range("G3").Value = 23.630375
$ws.Range("H3").Value = 70.89112499999999
$ws.Range("I3").Value = 0.9002398112414131
$ws.Range("J3").Value = 0.9002398112414129
$ws.Range("O3").Value = 0.004437346842596906
$ws.Range("P3").Value = 0.004437346842596906
$ws.Range("Q3").Value = 2.239482145916666
$ws.Range("R3").Value = 20.15533931324999
$ws.Range("S3").Value = 0.003994676283992119
$ws.Range("T3").Value = 0.003994676283992119

# --- Row 4: ECs -> MuSCs -------------------------------------------------
$ws.Range("G4").Value = 23.630375
$ws.Range("H4").Value = 70.89112499999999
$ws.Range("I4").Value = 0.9002398112414131
$ws.Range("J4").Value = 0.9002398112414129
$ws.Range("O4").Value = 0.4871601242408422
$ws.Range("P4").Value = 0.4871601242408422
$ws.Range("Q4").Value = 245.864576094625
$ws.Range("R4").Value = 2212.781184851624
$ws.Range("S4").Value = 0.4385609382909191
$ws.Range("T4").Value = 0.4385609382909191

# --- Row 5: FAPs -> ECs -------------------------------------------------
$ws.Range("I5").Value = 0.06214870537054815
$ws.Range("J5").Value = 0.06214870537054815
$ws.Range("M5").Value = 10.858287
$ws.Range("N5").Value = 32.574861
$ws.Range("O5").Value = 0.5084025289165609
$ws.Range("P5").Value = 0.508402528916561
$ws.Range("Q5").Value = 17.71355791458
$ws.Range("R5").Value = 159.42202123122
$ws.Range("S5").Value = 0.03159655897927693
$ws.Range("T5").Value = 0.03159655897927693

# --- Row 6: FAPs -> FAPs -------------------------------------------------
$ws.Range("I6").Value = 0.06214870537054815
$ws.Range("J6").Value = 0.06214870537054815
$ws.Range("O6").Value = 0.004437346842596906
$ws.Range("P6").Value = 0.004437346842596906
$ws.Range("S6").Value = 0.0002757753615474872
$ws.Range("T6").Value = 0.0002757753615474872

# --- Row 7: FAPs -> MuSCs -------------------------------------------------
$ws.Range("I7").Value = 0.06214870537054815
$ws.Range("J7").Value = 0.06214870537054815
$ws.Range("O7").Value = 0.4871601242408422
$ws.Range("P7").Value = 0.4871601242408422
$ws.Range("S7").Value = 0.03027637102972373
$ws.Range("T7").Value = 0.03027637102972373

# --- Row 8: MuSCs -> ECs -------------------------------------------------
$ws.Range("I8").Value = 0.03761148338803896
$ws.Range("J8").Value = 0.03761148338803896
$ws.Range("M8").Value = 10.858287
$ws.Range("N8").Value = 32.574861
$ws.Range("O8").Value = 0.5084025289165609
$ws.Range("P8").Value = 0.508402528916561
$ws.Range("Q8").Value = 10.719984998481
$ws.Range("R8").Value = 96.47986498632899
$ws.Range("S8").Value = 0.01912177327078223
$ws.Range("T8").Value = 0.01912177327078223

# --- Row 9: MuSCs -> FAPs -------------------------------------------------
$ws.Range("I9").Value = 0.03761148338803896
$ws.Range("J9").Value = 0.03761148338803896
$ws.Range("O9").Value = 0.004437346842596906
$ws.Range("P9").Value = 0.004437346842596906
$ws.Range("S9").Value = 0.0001668951970573007
$ws.Range("T9").Value = 0.0001668951970573007

# --- Row 10: MuSCs -> MuSCs -------------------------------------------------
$ws.Range("I10").Value = 0.03761148338803896
$ws.Range("J10").Value = 0.03761148338803896
$ws.Range("O10").Value = 0.4871601242408422
$ws.Range("P10").Value = 0.4871601242408422
$ws.Range("S10").Value = 0.01832281492019943
$ws.Range("T10").Value = 0.01832281492019943
